# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted at row 299 (shifting the
# existing rows 299-331 down to 300-332) and populated with the new
# observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 299, pushing all rows below it
# (299-331) down by one (to 300-332).
$ws.Rows.Item(299).Insert()

# Populate the newly inserted row 299 with the new data point.
$ws.Cells.Item(299, 1).Value  = 11
$ws.Cells.Item(299, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(299, 3).Value  = "Bíobío"
$ws.Cells.Item(299, 4).Value  = 45127
$ws.Cells.Item(299, 5).Value  = 8
$ws.Cells.Item(299, 6).Value  = 100112040
$ws.Cells.Item(299, 7).Value  = "Cilantro"
$ws.Cells.Item(299, 8).Value  = "Sin especificar"
$ws.Cells.Item(299, 9).Value  = "Primera"
$ws.Cells.Item(299, 10).Value = 90
$ws.Cells.Item(299, 11).Value = 7500
$ws.Cells.Item(299, 12).Value = 8000
$ws.Cells.Item(299, 13).Value = 7722
$ws.Cells.Item(299, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(299, 15).Value = "Región Metropolitana"
$ws.Cells.Item(299, 16).Value = 214
$ws.Cells.Item(299, 17).Value = 36
$ws.Cells.Item(299, 18).Value = "Hortaliza"
